# Generate Report for Archive
#
# 1) The localization status changes from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F3 and the
#    Status column, column C, on each language sheet).
# 2) Those status columns are narrowed to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Narrowed target column width (raw OOXML character width ~13.41,
# which is what ColumnWidth = 12.5 resolves to on save).
$newColWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newColWidth
